$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update the "Provide" value for both rows from Provide2018 -> Provide2012
$ws.Range("D2").Value = "Provide2012"
$ws.Range("D3").Value = "Provide2012"

# Update sheet view: scroll so column F is leftmost, and select R2
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("R2").Select()
